$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's text value while preserving its original cell style.
# Prices that look like plain numbers need a leading apostrophe so Excel keeps
# storing them as text (matching the existing text-formatted Price column),
# and re-applying the original Style afterwards avoids picking up the implicit
# "quote prefix" number format that Value-assignment would otherwise introduce.
function Set-TextValue($range, [string]$text) {
    $origStyle = $range.Style
    $range.Value = $text
    $range.Style = $origStyle
}

# Update Price (D) and Volume(1h) (E) columns for rows with changed values.
$ws.Range("D2").Value = "63.630.38"
$ws.Range("D3").Value = "2.614.19"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "'591.95"
$ws.Range("E5").Value = "  -1.60%  "
Set-TextValue $ws.Range("D6") "'150.13"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  +0.04%  "
Set-TextValue $ws.Range("D8") "'0.584"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("E9").Value = "  +0.57%  "
Set-TextValue $ws.Range("D10") "'5.78"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("E12").Value = "  +0.55%  "
Set-TextValue $ws.Range("D13") "'27.79"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "3.084.65"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "63.430.50"
$ws.Range("E15").Value = "  -0.87%  "
Set-TextValue $ws.Range("D16") "'0.0000159"
$ws.Range("E16").Value = "  +6.19%  "
$ws.Range("D17").Value = "2.603.38"
$ws.Range("E17").Value = "  -0.95%  "
Set-TextValue $ws.Range("D18") "'12.13"
$ws.Range("E18").Value = "  -0.67%  "
Set-TextValue $ws.Range("D19") "'4.76"
$ws.Range("E19").Value = "  +2.42%  "
Set-TextValue $ws.Range("D20") "'345.98"
$ws.Range("E20").Value = "  -1.32%  "
Set-TextValue $ws.Range("D21") "'6.93"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  +0.16%  "
Set-TextValue $ws.Range("D23") "'67.18"
$ws.Range("E23").Value = "  +1.10%  "
Set-TextValue $ws.Range("D24") "'1.69"
$ws.Range("E24").Value = "  -3.26%  "
Set-TextValue $ws.Range("D25") "'9.24"
$ws.Range("E25").Value = "  -0.13%  "
Set-TextValue $ws.Range("D26") "'1.66"
$ws.Range("E26").Value = "  -1.26%  "
Set-TextValue $ws.Range("D27") "'8.47"
$ws.Range("E27").Value = "  +3.43%  "
Set-TextValue $ws.Range("D28") "'547.78"
$ws.Range("E28").Value = "  +1.69%  "
Set-TextValue $ws.Range("D29") "'0.162"
$ws.Range("E29").Value = "  -1.46%  "
Set-TextValue $ws.Range("D30") "'1.00"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("E33").Value = "  +1.98%  "
Set-TextValue $ws.Range("D34") "'5.39"
$ws.Range("E34").Value = "  +1.67%  "
Set-TextValue $ws.Range("D35") "'6.12"
$ws.Range("E35").Value = "  -0.34%  "
Set-TextValue $ws.Range("D36") "'165.02"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E37").Value = "  +0.98%  "
Set-TextValue $ws.Range("D40") "'19.55"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").Value = "  -0.04%  "
Set-TextValue $ws.Range("D42") "'165.87"
$ws.Range("E42").Value = "  -1.84%  "
Set-TextValue $ws.Range("D43") "'4.08"
$ws.Range("E43").Value = "  +3.84%  "
Set-TextValue $ws.Range("D44") "'23.21"
$ws.Range("E44").Value = "  +7.61%  "
Set-TextValue $ws.Range("D45") "'0.0584"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  +7.33%  "
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  +17.98%  "

# Rows 38 and 39 swap coin identity (Coin/Link) and get new Price/Volume values
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D38") "'1.98"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D39") "'0.999"
$ws.Range("E39").Value = "  -0.07%  "
